$p = $ppt.ActivePresentation

# Slide 14 ("Grading Criteria") - Content Placeholder 2
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# 1) "Correctness" -> "correctness"
$full = $tr.Text
$idx = $full.IndexOf("Correctness")
$sub = $tr.Characters($idx + 1, "Correctness".Length)
$sub.Text = "correctness"

# 2) "Quality of design" -> "quality of design"
$full = $tr.Text
$idx = $full.IndexOf("Quality of design")
$sub = $tr.Characters($idx + 1, "Quality of design".Length)
$sub.Text = "quality of design"

# 3) "Adherence to our coding and documentation standards"
#    -> "adherence " + "to our coding and documentation standards" (two runs)
$full = $tr.Text
$idx = $full.IndexOf("Adherence ")
$sub = $tr.Characters($idx + 1, "Adherence ".Length)
$sub.Text = "adherence "
